$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.218.71'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '2.003.07'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''246.23'
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('D6').Value = '''0.630'
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('D7').Value = '''60.39'
$ws.Range('E7').Value = '  +3.85%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.385'
$ws.Range('E9').Value = '  +2.76%  '
$ws.Range('D10').Value = '''0.0806'
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').Value = '''15.12'
$ws.Range('E12').Value = '  +7.83%  '
$ws.Range('E13').Value = '  +7.73%  '
$ws.Range('D14').Value = '''0.851'
$ws.Range('E14').Value = '  +1.84%  '
$ws.Range('D15').Value = '2.295.95'
$ws.Range('E15').Value = '  +2.41%  '
$ws.Range('D16').Value = '''5.47'
$ws.Range('E16').Value = '  +3.69%  '
$ws.Range('D17').Value = '1.995.89'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = '37.139.19'
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('D19').Value = '''70.42'
$ws.Range('E19').Value = '  +1.03%  '
$ws.Range('D20').Value = '0.0₃0866'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').Value = '''5.20'
$ws.Range('E21').Value = '  +3.29%  '
$ws.Range('D22').Value = '''230.84'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '''2.47'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').Value = '''2.36'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').Value = '''9.47'
$ws.Range('E26').Value = '  +3.80%  '
$ws.Range('D27').Value = '''0.144'
$ws.Range('E27').Value = '  +3.75%  '
$ws.Range('D28').Value = '''164.26'
$ws.Range('E28').Value = '  +2.55%  '
$ws.Range('D29').Value = '''19.68'
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('E30').Value = '  +14.34%  '
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('E32').Value = '  +2.23%  '
$ws.Range('E34').Value = '  +3.36%  '
$ws.Range('D35').Value = '''2.39'
$ws.Range('E35').Value = '  +5.19%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  +2.70%  '
$ws.Range('D38').Value = '''3.31'
$ws.Range('E38').Value = '  -4.48%  '
$ws.Range('D39').Value = '''5.43'
$ws.Range('E39').Value = '  +2.58%  '
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').Value = '''2.91'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('D44').Value = '''16.88'
$ws.Range('E44').Value = '  +7.28%  '
$ws.Range('D45').Value = '''91.28'
$ws.Range('E45').Value = '  +4.00%  '
$ws.Range('D46').Value = '1.377.17'
$ws.Range('E46').Value = '  +0.42%  '
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('D48').Value = '''7.38'
$ws.Range('E48').Value = '  +3.74%  '
$ws.Range('D49').Value = '''2.05'
$ws.Range('E49').Value = '  +14.89%  '
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('D51').Value = '''46.51'
$ws.Range('E51').Value = '  +6.04%  '
